# 7rmartSupermarket/TestData.xlsx edit
#
# "url" sheet (tab order index 2): row 5 used to point at a local
# screenshot file ("ImageUrl" -> C:\...\code.png). Replace it with a
# real "ManageProduct" admin link, matching the other rows in that
# table, and wire up the hyperlink + selection like the existing rows.
#
# "ManageOrderPage" sheet (tab order index 3): becomes the active /
# selected sheet instead of "url", and its sample Order Id changes.

$wb = $excel.ActiveWorkbook

# --- "url" sheet: add the ManageProduct row -------------------------------
$wsUrl = $wb.Worksheets.Item("url")

# Set B5 (the URL) before A5 (the label) so the shared-string table grows
# in the same order as the target workbook (url string, then label string).
$wsUrl.Range("B5").Value = "https://groceryapp.uniqassosiates.com/admin/list-product"
$wsUrl.Range("A5").Value = "ManageProduct"

# Turn B5 into a real hyperlink, like B1:B4 above it.
$wsUrl.Hyperlinks.Add($wsUrl.Range("B5"), "https://groceryapp.uniqassosiates.com/admin/list-product")
$wsUrl.Range("B5").Style = "Hyperlink"

# Move the sheet's own selection onto the new label cell.
$wsUrl.Range("A5").Select()

# --- "ManageOrderPage" sheet: new sample value + becomes the active tab --
$wsOrder = $wb.Worksheets.Item("ManageOrderPage")
$wsOrder.Range("A2").Value = 133

# Activating this sheet makes it the workbook's active tab (and clears
# tabSelected/updates the selection on the previously active "url" sheet).
$wsOrder.Activate()
$wsOrder.Range("A2").Select()
